# Update cryptos list: refreshed prices and 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below are plain decimal-looking strings (e.g. "0.286"); force them to
# stay text (matching the sheet's existing text-formatted price column) instead
# of being auto-converted to numbers by Excel.
$textCells = @("D5", "D8", "D9", "D10", "D14", "D18", "D19", "D21", "D25", "D26", "D30", "D31", "D32", "D36", "D40", "D41", "D44", "D46", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "33.987.53"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.783.37"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "221.51"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "31.45"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "0.286"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "0.0709"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "2.039.72"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.778.53"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "10.51"
$ws.Range("E14").Value = "  -5.12%  "
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "33.975.55"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "67.92"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "244.62"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").Value = "157.53"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Value = "16.38"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "0.0521"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.68"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "1.402.47"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").Value = "0.638"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("D40").Value = "79.33"
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "5.93"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.03"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.938.96"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("E51").Value = "  -0.89%  "
